$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("C1").Value = "variable_trajectory_group"
    $ws.Range("D1").Value = "normalize_group"
    $ws.Range("E1").Value = "trajgroup_no_vary_q"
    $ws.Range("F1").Value = "uniform_scaling_q"
}
